# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (column I) and DialogAct (column J) values for the rows whose
# annotations changed after the transcript clean up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 8;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 20; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 34; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 35; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 37; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 51; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 55; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 74; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 75; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 76; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 80; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 81; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.Tag
    $ws.Range("J$($u.Row)").Value = $u.Act
}
